# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" rows (B16:J19) are re-sorted by
# "Periodo Mora" ascending (1908, 2009, 2112, 2206), each period keeping
# its paired "Valor Mora" amount.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Capture the current Periodo Mora / Valor Mora pairs (rows 16..19)
$rows = 16..19
$pairs = @()
foreach ($r in $rows) {
    $periodo = $ws.Range("E$r").Value()
    $valor = $ws.Range("F$r").Value()
    $pairs += , @($periodo, $valor)
}

# Sort ascending by Periodo Mora (stored as text, so sort numerically)
$sorted = $pairs | Sort-Object { [int]$_[0] }

# Write the sorted pairs back into the same rows, keeping everything else intact
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Range("E$r").Value = [string]$sorted[$i][0]
    $ws.Range("F$r").Value = $sorted[$i][1]
}
